$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3375
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3375
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 3375
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -4005
$ws.Range("H79").Value = 3375
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3375
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 3375
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -5559
$ws.Range("H105").Value = 71396
$ws.Range("J105").Value = 71396
$ws.Range("L105").Value = 71396
$ws.Range("N105").Value = -78384
$ws.Range("H116").Value = 2406280.5
$ws.Range("I116").Value = 10991054
$ws.Range("J116").Value = 2543.84
$ws.Range("K116").Value = 10991054
$ws.Range("L116").Value = 2543.84
$ws.Range("M116").Value = -10987612
$ws.Range("N116").Value = -9427.84
$ws.Range("H132").Value = 2456.38
$ws.Range("I132").Value = 1648.1086
$ws.Range("K132").Value = 4944.325800000001
$ws.Range("M132").Value = -2414.325800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1266.1471
$ws.Range("I2").Value = 1256.4138
$ws.Range("J2").Value = 1322.6
$ws.Range("K2").Value = 1256.4138
$ws.Range("L2").Value = 1322.6
$ws.Range("M2").Value = -1143.4138
$ws.Range("N2").Value = -1548.6
$ws.Range("H63").Value = 1089.2858
$ws.Range("I63").Value = 1089.2858
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1089.2858
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -403.2858000000001
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1089.2858
$ws.Range("I66").Value = 1089.2858
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 5446.429
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -2014.429
$ws.Range("N66").ClearContents()
$ws.Range("H101").Value = 53333.332
$ws.Range("J101").Value = 53333.332
$ws.Range("L101").Value = 53333.332
$ws.Range("N101").Value = -59823.332
$ws.Range("H116").Value = 1266.1471
$ws.Range("I116").Value = 1256.4138
$ws.Range("J116").Value = 1322.6
$ws.Range("K116").Value = 1256.4138
$ws.Range("L116").Value = 1322.6
$ws.Range("M116").Value = 1037.5862
$ws.Range("N116").Value = -5910.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1266.1471
$ws.Range("I3").Value = 1256.4138
$ws.Range("J3").Value = 1322.6
$ws.Range("K3").Value = 1256.4138
$ws.Range("L3").Value = 1322.6
$ws.Range("M3").Value = -1142.4138
$ws.Range("N3").Value = -1550.6
$ws.Range("H35").Value = 17018.5
$ws.Range("J35").Value = 17018.5
$ws.Range("L35").Value = 17018.5
$ws.Range("N35").Value = -17638.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 10571.429
$ws.Range("J50").Value = 10571.429
$ws.Range("L50").Value = 10571.429
$ws.Range("N50").Value = -11821.429
$ws.Range("H51").Value = 10857.143
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 10857.143
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 10857.143
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -12329.143
$ws.Range("H59").Value = 15034.692
$ws.Range("J59").Value = 15034.692
$ws.Range("L59").Value = 15034.692
$ws.Range("N59").Value = -17324.692
$ws.Range("H60").Value = 10000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 10000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 10000
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -11022
$ws.Range("H61").Value = 10857.143
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 10857.143
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 10857.143
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -11553.143
$ws.Range("H68").Value = 34950
$ws.Range("J68").Value = 34950
$ws.Range("L68").Value = 34950
$ws.Range("N68").Value = -36448
$ws.Range("H71").Value = 34950
$ws.Range("J71").Value = 34950
$ws.Range("L71").Value = 104850
$ws.Range("N71").Value = -112338
$ws.Range("H122").Value = 1801
$ws.Range("I122").Value = 1361.4286
$ws.Range("J122").Value = 2826.6667
$ws.Range("K122").Value = 4084.2858
$ws.Range("L122").Value = 8480.000100000001
$ws.Range("M122").Value = -1634.2858
$ws.Range("N122").Value = -13380.0001
$ws.Range("H132").Value = 1830.762
$ws.Range("I132").Value = 1496.2333
$ws.Range("J132").Value = 2667.0833
$ws.Range("K132").Value = 4488.699900000001
$ws.Range("L132").Value = 8001.249899999999
$ws.Range("M132").Value = -1958.699900000001
$ws.Range("N132").Value = -13061.2499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 4143
$ws.Range("J105").Value = 4143
$ws.Range("L105").Value = 12429
$ws.Range("N105").Value = -17671
$ws.Range("H107").Value = 1250435.6
$ws.Range("I107").Value = 460.2
$ws.Range("J107").Value = 3333728
$ws.Range("K107").Value = 1380.6
$ws.Range("L107").Value = 10001184
$ws.Range("M107").Value = 539.4000000000001
$ws.Range("N107").Value = -10005024
$ws.Range("H109").Value = 4265.654
$ws.Range("I109").Value = 602.3333
$ws.Range("J109").Value = 4743.478
$ws.Range("K109").Value = 1806.9999
$ws.Range("L109").Value = 14230.434
$ws.Range("M109").Value = -766.9999
$ws.Range("N109").Value = -16310.434
$ws.Range("H131").Value = 748.4925500000001
$ws.Range("I131").Value = 370.3889
$ws.Range("J131").Value = 887.38776
$ws.Range("K131").Value = 1111.1667
$ws.Range("L131").Value = 2662.16328
$ws.Range("M131").Value = 3928.8333
$ws.Range("N131").Value = -12742.16328

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2424
$ws.Range("I31").Value = 2424
$ws.Range("K31").Value = 2424
$ws.Range("M31").Value = -2132
$ws.Range("H37").Value = 2424
$ws.Range("I37").Value = 2424
$ws.Range("K37").Value = 2424
$ws.Range("M37").Value = -2147
$ws.Range("H118").Value = 10100
$ws.Range("J118").Value = 10100
$ws.Range("L118").Value = 10100
$ws.Range("N118").Value = -13414

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1972.25
$ws.Range("I61").Value = 1972.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1972.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1770.25
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 1972.25
$ws.Range("I113").Value = 1972.25
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1972.25
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 197.75
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3290.973
$ws.Range("I122").Value = 3236.3635
$ws.Range("J122").Value = 3371.0667
$ws.Range("K122").Value = 9709.0905
$ws.Range("L122").Value = 10113.2001
$ws.Range("M122").Value = -7259.0905
$ws.Range("N122").Value = -15013.2001
$ws.Range("H128").Value = 31867.285
$ws.Range("J128").Value = 31867.285
$ws.Range("L128").Value = 31867.285
$ws.Range("N128").Value = -41827.285
$ws.Range("H136").Value = 3214.3438
$ws.Range("I136").Value = 1889.6
$ws.Range("J136").Value = 4383.2354
$ws.Range("K136").Value = 5668.799999999999
$ws.Range("L136").Value = 13149.7062
$ws.Range("M136").Value = -3118.799999999999
$ws.Range("N136").Value = -18249.7062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 304.79166
$ws.Range("I113").Value = 284.77274
$ws.Range("J113").Value = 525
$ws.Range("K113").Value = 854.31822
$ws.Range("L113").Value = 1575
$ws.Range("M113").Value = 1315.68178
$ws.Range("N113").Value = -5915
$ws.Range("H132").Value = 1652.4918
$ws.Range("I132").Value = 812.2368
$ws.Range("J132").Value = 3040.739
$ws.Range("K132").Value = 2436.7104
$ws.Range("L132").Value = 9122.217000000001
$ws.Range("M132").Value = 93.28960000000006
$ws.Range("N132").Value = -14182.217
$ws.Range("H136").Value = 5236.037
$ws.Range("I136").Value = 6249.8887
$ws.Range("K136").Value = 18749.6661
$ws.Range("M136").Value = -16199.6661
